$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix a stray description value in row 5 (was "20", should be "facilities")
$ws.Range("I5").Value = "facilities"

# Row 6: The Wyndham New Yorker Hotel (id 105)
$ws.Range("A6:H6").NumberFormat = "@"
$ws.Range("B6").Value = "The Wyndham New Yorker Hotel"
$ws.Range("C6").Value = "United States"
$ws.Range("D6").Value = "New York"
$ws.Range("E6").Value = "481 Eighth Avenue"
$ws.Range("F6").Value = "8888881"
$ws.Range("G6").Value = "5"
$ws.Range("H6").Value = "300"
$ws.Range("I6").Value = "Fitness"

# Row 7: Hotel Mela Times Square (id 106)
$ws.Range("B7").Value = "Hotel Mela Times Square"
$ws.Range("C7").Value = "United States"
$ws.Range("D7").Value = "New York"
$ws.Range("E7").Value = "120 West 44th Street"
$ws.Range("F7").Value = "1829981"
$ws.Range("G7").Value = "5"
$ws.Range("H7").Value = "600"
$ws.Range("I7").Value = "everything"

# Row 8: Hotel Pennsylvania (id 107)
$ws.Range("B8").Value = "Hotel Pennsylvania"
$ws.Range("C8").Value = "United States"
$ws.Range("D8").Value = "New York"
$ws.Range("E8").Value = "401 7th Avenue"
$ws.Range("F8").Value = "999999"
$ws.Range("G8").Value = "5"
$ws.Range("H8").Value = "400"
$ws.Range("I8").Value = "everything"

# Row 9: The Savoy Hotel (id 108)
$ws.Range("B9").Value = "The Savoy Hotel"
$ws.Range("C9").Value = "United Kingdom"
$ws.Range("D9").Value = "London"
$ws.Range("E9").Value = "Strand, West End Soho"
$ws.Range("F9").Value = "11233"
$ws.Range("G9").Value = "3"
$ws.Range("H9").Value = "40"
$ws.Range("I9").Value = "---"

# Row 10: Days Inn Hilton Head (id 109) - new row
$ws.Range("A10").Value = "109"
$ws.Range("B10").Value = "Days Inn Hilton Head"
$ws.Range("C10").Value = "United States"
$ws.Range("D10").Value = "hilton"
$ws.Range("E10").Value = "9 Marina Side Drive"
$ws.Range("F10").Value = "999999"
$ws.Range("G10").Value = "5"
$ws.Range("H10").Value = "200"
$ws.Range("I10").Value = "fitness"

# Row 11: Hilton Head Marriott Resort & Spa (id 110) - new row
$ws.Range("A11").Value = "110"
$ws.Range("B11").Value = "Hilton Head Marriott Resort & Spa"
$ws.Range("C11").Value = "United States"
$ws.Range("D11").Value = "Hilton"
$ws.Range("E11").Value = "---"
$ws.Range("F11").Value = "11111"
$ws.Range("G11").Value = "5"
$ws.Range("H11").Value = "500"
$ws.Range("I11").Value = "unknown"

$ws.Range("I10").Select()
